$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 601 (existing rows 601.. shift down to 603..667)
$ws.Rows.Item(601).Insert()
$ws.Rows.Item(601).Insert()

# --- New row 601: Early Majestic / Especial ---
$ws.Cells.Item(601, 1).Value = 5
$ws.Cells.Item(601, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(601, 3).Value = "Maule"
$ws.Cells.Item(601, 4).Value = 45265
$ws.Cells.Item(601, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(601, 5).Value = 7
$ws.Cells.Item(601, 6).Value = "Fruta"
$ws.Cells.Item(601, 7).Value = 100103
$ws.Cells.Item(601, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(601, 9).Value = 100103004
$ws.Cells.Item(601, 10).Value = "Durazno"
$ws.Cells.Item(601, 11).Value = "Early Majestic"
$ws.Cells.Item(601, 12).Value = "Especial"
$ws.Cells.Item(601, 13).Value = 300
$ws.Cells.Item(601, 14).Value = 15000
$ws.Cells.Item(601, 15).Value = 15000
$ws.Cells.Item(601, 16).Value = 15000
$ws.Cells.Item(601, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(601, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(601, 19).Value = 1000
$ws.Cells.Item(601, 20).Value = 15

# --- New row 602: Early Treat / Especial ---
$ws.Cells.Item(602, 1).Value = 5
$ws.Cells.Item(602, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(602, 3).Value = "Maule"
$ws.Cells.Item(602, 4).Value = 45265
$ws.Cells.Item(602, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(602, 5).Value = 7
$ws.Cells.Item(602, 6).Value = "Fruta"
$ws.Cells.Item(602, 7).Value = 100103
$ws.Cells.Item(602, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(602, 9).Value = 100103004
$ws.Cells.Item(602, 10).Value = "Durazno"
$ws.Cells.Item(602, 11).Value = "Early Treat"
$ws.Cells.Item(602, 12).Value = "Especial"
$ws.Cells.Item(602, 13).Value = 250
$ws.Cells.Item(602, 14).Value = 15000
$ws.Cells.Item(602, 15).Value = 15000
$ws.Cells.Item(602, 16).Value = 15000
$ws.Cells.Item(602, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(602, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(602, 19).Value = 1000
$ws.Cells.Item(602, 20).Value = 15

Write-Output "done"
